$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "0.1225487848221539"
$ws.Range("B3").Value = "0.0006642914287983296"
$ws.Range("C3").Value = "0.000701332145311826"
$ws.Range("D3").Value = "4.008249250378249"
$ws.Range("E3").Value = "0.1513892964512931"
$ws.Range("F3").Value = "-0.0007102998848576619"
$ws.Range("G3").Value = "0.002038882742454322"
$ws.Range("H3").Value = "0.1232130762509522"
$ws.Range("B4").Value = "0.007231600581606429"
$ws.Range("C4").Value = "0.0008914135766871937"
$ws.Range("D4").Value = "5.958503751565218"
$ws.Range("E4").Value = "0.02317579185347605"
$ws.Range("F4").Value = "0.005484455147701155"
$ws.Range("G4").Value = "0.008978746015511704"
$ws.Range("H4").Value = "0.1297803854037603"
$ws.Range("B5").Value = "0.01774536296681863"
$ws.Range("C5").Value = "0.001642302739317422"
$ws.Range("D5").Value = "12.63977650596076"
$ws.Range("E5").Value = "0.002906319831651242"
$ws.Range("F5").Value = "0.01452649711695079"
$ws.Range("G5").Value = "0.02096422881668646"
$ws.Range("H5").Value = "0.1402941477889725"
$ws.Range("B6").Value = "0.03937598183912521"
$ws.Range("C6").Value = "0.002254146770695159"
$ws.Range("D6").Value = "19.18421404728755"
$ws.Range("E6").Value = "0.04240242214453866"
$ws.Range("F6").Value = "0.03495791795192704"
$ws.Range("G6").Value = "0.04379404572632338"
$ws.Range("H6").Value = "0.1619247666612791"
$ws.Range("B7").Value = "0.04944733184074588"
$ws.Range("C7").Value = "0.002575337438793602"
$ws.Range("D7").Value = "18.47603540663996"
$ws.Range("E7").Value = "0.04749048197975473"
$ws.Range("F7").Value = "0.04439974714467964"
$ws.Range("G7").Value = "0.05449491653681209"
$ws.Range("H7").Value = "0.1719961166628998"
$ws.Range("B8").Value = "0.05038394550115362"
$ws.Range("C8").Value = "0.003838852250574208"
$ws.Range("D8").Value = "15.74824721934973"
$ws.Range("E8").Value = "0.1208775797372662"
$ws.Range("F8").Value = "0.04285990862410705"
$ws.Range("G8").Value = "0.05790798237820018"
$ws.Range("H8").Value = "0.1729327303233075"
$ws.Range("B9").Value = "0.04935570879665811"
$ws.Range("C9").Value = "0.005108216307462673"
$ws.Range("D9").Value = "14.47308383289499"
$ws.Range("E9").Value = "0.1373486666019877"
$ws.Range("F9").Value = "0.03934376161083224"
$ws.Range("G9").Value = "0.05936765598248396"
$ws.Range("H9").Value = "0.171904493618812"
$ws.Range("B10").Value = "-0.1225487848221539"
$ws.Range("C10").Value = "0.0005637990335111847"
$ws.Range("D10").Value = "-244.8194557402864"
$ws.Range("E10").Value = "0"
$ws.Range("F10").Value = "-0.1236538151594107"
$ws.Range("G10").Value = "-0.1214437544848971"
$ws.Range("B11").Value = "-0.06057886980565406"
$ws.Range("C11").Value = "0.0006286346230349478"
$ws.Range("D11").Value = "-100.3116957944558"
$ws.Range("E11").Value = "0"
$ws.Range("F11").Value = "-0.06181097596994985"
$ws.Range("G11").Value = "-0.0593467636413583"
$ws.Range("H11").Value = "0.06196991501649983"
$ws.Range("B12").Value = "-0.05231118864881306"
$ws.Range("C12").Value = "0.0006277323703255119"
$ws.Range("D12").Value = "-86.74497594235785"
$ws.Range("E12").Value = "0"
$ws.Range("F12").Value = "-0.05354152644009074"
$ws.Range("G12").Value = "-0.05108085085753535"
$ws.Range("H12").Value = "0.07023759617334083"
$ws.Range("B13").Value = "-0.04361797780264658"
$ws.Range("C13").Value = "0.0006194816333964756"
$ws.Range("D13").Value = "-71.43910871122776"
$ws.Range("E13").Value = "0"
$ws.Range("F13").Value = "-0.04483214439492577"
$ws.Range("G13").Value = "-0.04240381121036739"
$ws.Range("H13").Value = "0.07893080701950732"
$ws.Range("B14").Value = "-0.03949396645757155"
$ws.Range("C14").Value = "0.0006173284884751704"
$ws.Range("D14").Value = "-63.8633445034966"
$ws.Range("E14").Value = "5.723627552016212e-261"
$ws.Range("F14").Value = "-0.04070391295890632"
$ws.Range("G14").Value = "-0.03828401995623677"
$ws.Range("H14").Value = "0.08305481836458234"
$ws.Range("B15").Value = "-0.03464392156095113"
$ws.Range("C15").Value = "0.000601832964539926"
$ws.Range("D15").Value = "-56.38244097199102"
$ws.Range("E15").Value = "3.488455357346417e-140"
$ws.Range("F15").Value = "-0.03582349728424131"
$ws.Range("G15").Value = "-0.03346434583766096"
$ws.Range("H15").Value = "0.08790486326120277"
$ws.Range("B16").Value = "-0.03189213366479873"
$ws.Range("C16").Value = "0.0005806754641529313"
$ws.Range("D16").Value = "-52.28810489118829"
$ws.Range("E16").Value = "1.628779986981933e-113"
$ws.Range("F16").Value = "-0.03303024134113088"
$ws.Range("G16").Value = "-0.03075402598846658"
$ws.Range("H16").Value = "0.09065665115735516"
$ws.Range("B17").Value = "-0.02900455138610384"
$ws.Range("C17").Value = "0.0005856438183715592"
$ws.Range("D17").Value = "-47.49473851410416"
$ws.Range("E17").Value = "4.151123576442654e-71"
$ws.Range("F17").Value = "-0.03015239688676554"
$ws.Range("G17").Value = "-0.02785670588544216"
$ws.Range("H17").Value = "0.09354423343605005"
$ws.Range("B18").Value = "-0.02778046200621384"
$ws.Range("C18").Value = "0.0005918602127313485"
$ws.Range("D18").Value = "-45.36702410947227"
$ws.Range("E18").Value = "1.291535726137029e-44"
$ws.Range("F18").Value = "-0.02894049143467792"
$ws.Range("G18").Value = "-0.02662043257774977"
$ws.Range("H18").Value = "0.09476832281594005"
$ws.Range("B19").Value = "-0.02266152359602177"
$ws.Range("C19").Value = "0.0005999684431477382"
$ws.Range("D19").Value = "-37.2420782396519"
$ws.Range("E19").Value = "1.327004240520837e-35"
$ws.Range("F19").Value = "-0.02383744499539763"
$ws.Range("G19").Value = "-0.02148560219664589"
$ws.Range("H19").Value = "0.09988726122613213"
$ws.Range("B20").Value = "-0.01943808609083728"
$ws.Range("C20").Value = "0.0006047593180176778"
$ws.Range("D20").Value = "-30.70796846209864"
$ws.Range("E20").Value = "4.956202391358433e-16"
$ws.Range("F20").Value = "-0.02062339754251966"
$ws.Range("G20").Value = "-0.0182527746391549"
$ws.Range("H20").Value = "0.1031106987313166"
$ws.Range("B21").Value = "-0.01535763598063881"
$ws.Range("C21").Value = "0.0006064099760532624"
$ws.Range("D21").Value = "-24.77052937055897"
$ws.Range("E21").Value = "4.595342966176211e-08"
$ws.Range("F21").Value = "-0.01654618261416826"
$ws.Range("G21").Value = "-0.01416908934710935"
$ws.Range("H21").Value = "0.1071911488415151"
$ws.Range("B22").Value = "-0.01168270091330123"
$ws.Range("C22").Value = "0.0005915934827280048"
$ws.Range("D22").Value = "-18.9907449089927"
$ws.Range("E22").Value = "0.0004712429164848759"
$ws.Range("F22").Value = "-0.01284220758429973"
$ws.Range("G22").Value = "-0.01052319424230274"
$ws.Range("H22").Value = "0.1108660839088527"
$ws.Range("B23").Value = "-0.009026339830258931"
$ws.Range("C23").Value = "0.0005838807018699947"
$ws.Range("D23").Value = "-15.42972490255434"
$ws.Range("E23").Value = "0.05298309404894173"
$ws.Range("F23").Value = "-0.01017072968032115"
$ws.Range("G23").Value = "-0.007881949980196725"
$ws.Range("H23").Value = "0.113522444991895"
$ws.Range("B24").Value = "-0.006258683418861363"
$ws.Range("C24").Value = "0.0005709994668426012"
$ws.Range("D24").Value = "-10.05139853199129"
$ws.Range("E24").Value = "0.0001011949902866033"
$ws.Range("F24").Value = "-0.007377826487897082"
$ws.Range("G24").Value = "-0.005139540349825643"
$ws.Range("H24").Value = "0.1162901014032925"
$ws.Range("B25").Value = "-0.006472555300048615"
$ws.Range("C25").Value = "0.0005547277163917625"
$ws.Range("D25").Value = "-9.478052682572638"
$ws.Range("E25").Value = "0.03125499100556205"
$ws.Range("F25").Value = "-0.007559806294410761"
$ws.Range("G25").Value = "-0.005385304305686471"
$ws.Range("H25").Value = "0.1160762295221053"
$ws.Range("B26").Value = "0.06177481790951234"
$ws.Range("C26").Value = "0.003779857791944061"
$ws.Range("D26").Value = "46.02572417280035"
$ws.Range("E26").Value = "0.1414671711601413"
$ws.Range("F26").Value = "0.05436641309928673"
$ws.Range("G26").Value = "0.06918322271973795"
$ws.Range("H26").Value = "0.1843236027316662"
